$d = $word.ActiveDocument

$d.Content.Find.Execute(" 1,000", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 10,000", 2)

$d.Content.Find.Execute("Varies, Usually around 210-ish", $true, $false, $false, $false, $false,
                         $true, 1, $false, "219", 2)
